$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Price (D) and Volume(1h) (E) columns.
# D-column values that look numeric must be forced to text
# (NumberFormat "@") so Excel stores them exactly as the original
# author-supplied strings (e.g. "1.000" must stay "1.000", not 1).

$ws.Cells.Item(2, 4).Value = "23.405.88"
$ws.Cells.Item(2, 5).Value = "  -0.12%  "

$ws.Cells.Item(3, 4).Value = "1.626.39"
$ws.Cells.Item(3, 5).Value = "  -0.37%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9996"
$ws.Cells.Item(4, 5).Value = "  -0.45%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "1.000"
$ws.Cells.Item(5, 5).Value = "  -0.24%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "304.66"
$ws.Cells.Item(6, 5).Value = "  -1.28%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3787"
$ws.Cells.Item(7, 5).Value = "  +0.33%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "51.96"
$ws.Cells.Item(8, 5).Value = "  -1.81%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.3629"
$ws.Cells.Item(9, 5).Value = "  -1.04%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.229"
$ws.Cells.Item(10, 5).Value = "  -3.28%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.08099"
$ws.Cells.Item(11, 5).Value = "  -1.06%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.001"
$ws.Cells.Item(12, 5).Value = "  -0.33%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "22.69"
$ws.Cells.Item(13, 5).Value = "  -1.58%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.541"
$ws.Cells.Item(14, 5).Value = "  -1.60%  "

$ws.Cells.Item(15, 5).Value = "  -2.48%  "

$ws.Cells.Item(16, 5).Value = "  -3.18%  "

$ws.Cells.Item(17, 4).Value = "1.624.31"
$ws.Cells.Item(17, 5).Value = "  -0.90%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "93.50"
$ws.Cells.Item(18, 5).Value = "  -1.37%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06906"
$ws.Cells.Item(19, 5).Value = "  -0.67%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "17.88"
$ws.Cells.Item(20, 5).Value = "  -2.15%  "

$ws.Cells.Item(21, 5).Value = "  -0.23%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.412"
$ws.Cells.Item(22, 5).Value = "  -2.37%  "

$ws.Cells.Item(23, 4).Value = "23.398.30"
$ws.Cells.Item(23, 5).Value = "  -0.20%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "12.73"
$ws.Cells.Item(24, 5).Value = "  -1.85%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "3.235"
$ws.Cells.Item(25, 5).Value = "  +3.45%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.444"
$ws.Cells.Item(26, 5).Value = "  +1.33%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "21.13"
$ws.Cells.Item(27, 5).Value = "  -1.11%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "149.11"
$ws.Cells.Item(28, 5).Value = "  -1.43%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.288"
$ws.Cells.Item(29, 5).Value = "  -1.07%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "134.25"
$ws.Cells.Item(30, 5).Value = "  -1.26%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "2.306"
$ws.Cells.Item(31, 5).Value = "  -5.09%  "

$ws.Cells.Item(32, 4).Value = "1.809.14"
$ws.Cells.Item(32, 5).Value = "  -0.50%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "6.786"
$ws.Cells.Item(33, 5).Value = "  -0.07%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "11.05"
$ws.Cells.Item(34, 5).Value = "  +5.84%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.9515"
$ws.Cells.Item(35, 5).Value = "  -2.15%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.02793"
$ws.Cells.Item(36, 5).Value = "  -0.67%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.2529"
$ws.Cells.Item(37, 5).Value = "  -0.24%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.08838"
$ws.Cells.Item(38, 5).Value = "  -0.08%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "6.107"
$ws.Cells.Item(39, 5).Value = "  -1.73%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.07180"
$ws.Cells.Item(40, 5).Value = "  -3.30%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.359"
$ws.Cells.Item(41, 5).Value = "  -2.74%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.7065"
$ws.Cells.Item(42, 5).Value = "  -1.13%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "16.22"
$ws.Cells.Item(43, 5).Value = "  -0.14%  "

$ws.Cells.Item(44, 5).Value = "  -2.12%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.6460"
$ws.Cells.Item(45, 5).Value = "  -2.29%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.323"
$ws.Cells.Item(46, 5).Value = "  -1.43%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.9994"
$ws.Cells.Item(47, 5).Value = "  -0.07%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.998"
$ws.Cells.Item(48, 5).Value = "  -1.29%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.07994"
$ws.Cells.Item(49, 5).Value = "  -0.63%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.202"
$ws.Cells.Item(50, 5).Value = "  -1.13%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "125.73"
$ws.Cells.Item(51, 5).Value = "  -4.05%  "
